$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily format the Price column as Text so numeric-looking strings
# (e.g. "1.00", "64.20") are preserved exactly as strings rather than
# being auto-converted to numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "58.319.98"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "2.484.01"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").Value = "522.32"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").Value = "132.76"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "0.556"
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "2.515.19"
$ws.Range("E9").Value = "  +2.27%  "
$ws.Range("D10").Value = "0.0976"
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "5.13"
$ws.Range("E12").Value = "  -2.20%  "
$ws.Range("D13").Value = "0.332"
$ws.Range("E13").Value = "  -1.81%  "
$ws.Range("D14").Value = "2.941.69"
$ws.Range("E14").Value = "  +1.78%  "
$ws.Range("D15").Value = "58.289.08"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").Value = "22.16"
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "2.516.88"
$ws.Range("E18").Value = "  +2.59%  "
$ws.Range("D19").Value = "10.68"
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("D20").Value = "321.61"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("D22").Value = "6.04"
$ws.Range("E22").Value = "  +5.77%  "
$ws.Range("D23").Value = "0.996"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("D24").Value = "64.20"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").Value = "0.401"
$ws.Range("E25").Value = "  -1.60%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "0.161"
$ws.Range("E27").Value = "  +1.23%  "
$ws.Range("D29").Value = "0.0₃0754"
$ws.Range("E29").Value = "  +2.23%  "
$ws.Range("D30").Value = "1.72"
$ws.Range("E30").Value = "  +2.23%  "
$ws.Range("E31").Value = "  +3.41%  "
$ws.Range("D32").Value = "167.19"
$ws.Range("E32").Value = "  +0.41%  "
$ws.Range("D33").Value = "6.26"
$ws.Range("E33").Value = "  +1.29%  "
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "0.992"
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("D36").Value = "18.11"
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("E37").Value = "  -2.24%  "
$ws.Range("E38").Value = "  +0.60%  "
$ws.Range("D39").Value = "36.83"
$ws.Range("E39").Value = "  +1.89%  "
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("D41").Value = "0.781"
$ws.Range("E41").Value = "  -0.60%  "
$ws.Range("D42").Value = "279.06"
$ws.Range("E42").Value = "  +3.52%  "
$ws.Range("D43").Value = "3.50"
$ws.Range("E43").Value = "  +2.13%  "
$ws.Range("D44").Value = "5.08"
$ws.Range("E44").Value = "  +2.71%  "
$ws.Range("D45").Value = "0.598"
$ws.Range("E45").Value = "  +1.84%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "123.37"
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "0.0918"
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("D48").Value = "0.0501"
$ws.Range("E48").Value = "  +3.53%  "
$ws.Range("D49").Value = "17.80"
$ws.Range("E49").Value = "  +1.40%  "
$ws.Range("E50").Value = "  +1.86%  "
$ws.Range("D51").Value = "16.96"
$ws.Range("E51").Value = "  +1.76%  "

# Restore the default (unstyled) cell style now that values are set as text.
$priceRange.Style = "Normal"
